$wb = $excel.ActiveWorkbook

# --- DegreeRequirement sheet: convert the per-row CONCATENATE formula into a
#     shared formula (D2:D13), matching the commit's formula consolidation. ---
$wsReq = $wb.Worksheets.Item("DegreeRequirement")
$reqFormula = '=CONCATENATE("new DegreeRequirement{","DegreeRequirementID=",A:A,",","DegreeID=",B:B,",","RequirementID=",C:C,"},")'
$wsReq.Range("D2:D13").Formula = $reqFormula

# Move the selection off the old D2:D13 block onto D18 (matches the diff's
# saved selection state) and leave this sheet NOT the active tab.
$wsReq.Activate()
$wsReq.Range("D18").Select()

# --- DegreePlan sheet: rename the header labels (drop the "(PK)/(FK)/(U,n)"
#     annotations) and add a new column F with the DegreePlan{} builder
#     formula, mirroring the sibling Degree/DegreeRequirement/Student sheets. ---
$wsPlan = $wb.Worksheets.Item("DegreePlan")
$wsPlan.Range("A1").Value = "DegreePlanID"
$wsPlan.Range("B1").Value = "DegreeID"
$wsPlan.Range("C1").Value = "StudentID"
$wsPlan.Range("D1").Value = "DegreePlanAbbrev"
$wsPlan.Range("E1").Value = "DegreePlanName"

$planFormula = "=CONCATENATE(`"new DegreePlan{`",`"DegreePlanID=`",A:A,`",`",`"DegreeID=`",B:B,`",`",`"StudentID=`",C:C,`",`",`"DegreePlanAbbrev=`",`"''`",D:D,`"''`",`",`",`"DegreePlanName=`",`"''`",E:E,`"''`",`"},`")"
$wsPlan.Range("F2:F11").Formula = $planFormula

# Widen the new column to fit the generated C# snippet text.
$wsPlan.Columns.Item(6).ColumnWidth = 131

# This sheet becomes the active tab, with F2:F11 selected (matches activeTab="4").
$wsPlan.Activate()
$wsPlan.Range("F2:F11").Select()
